$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.434.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.84%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.641.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.36%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "

# Row 5
$ws.Range("E5").Value = "  -0.05%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "305.07"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.57%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3736"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.77%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "52.32"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.81%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3631"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.253"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.83%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08123"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.03%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9995"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.10%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.599"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.09%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001271"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.01%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.289"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.71%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.632.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.66%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.41%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06900"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.47%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.13%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.510"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.29%  "

# Row 22
$ws.Range("E22").Value = "  -0.01%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "23.445.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.89%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.30%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.088"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.67%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.413"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.67%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.11%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.27%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.335"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.42%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.84%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.278"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.07%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.811.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.63%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.827"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.52%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9515"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.04%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02811"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.13%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.85%  "

# Row 37
$ws.Range("E37").Value = "  -0.59%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.07236"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.51%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.111"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.07%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.08773"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.41%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.372"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.71%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7065"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.64%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.57%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6526"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.19%  "

# Row 46
$ws.Range("E46").Value = "  +0.53%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9992"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.07%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.010"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.12%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07970"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.25%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "128.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.07%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.200"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.41%  "

